$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Date column (BF): values were originally "1-23-2013-14" and should be
# normalized to ISO format "2014-01-23". Force Text format first so Excel keeps
# the literal string instead of re-interpreting it as a date serial number.
$ws.Range("BF2:BF31").NumberFormat = "@"
for ($r = 2; $r -le 31; $r++) {
    $ws.Range("BF$r").Value = "2014-01-23"
}

# Corrected per-team statistics (one day of NBA data had been shifted).
$ws.Range("AM2").Value = 4
$ws.Range("AT3").Value = 20
$ws.Range("AW3").Value = 24
$ws.Range("AD4").Value = 28
$ws.Range("AO5").Value = 9
$ws.Range("AO6").Value = 12
$ws.Range("BC6").Value = 15
$ws.Range("AD7").Value = 10
$ws.Range("AF7").Value = 24
$ws.Range("AV7").Value = 13
$ws.Range("AW7").Value = 25
$ws.Range("D9").Value = 40
$ws.Range("F9").Value = 20
$ws.Range("G9").Value = 0.5
$ws.Range("J9").Value = 84.90000000000001
$ws.Range("K9").Value = 0.45
$ws.Range("N9").Value = 0.36
$ws.Range("O9").Value = 18.8
$ws.Range("P9").Value = 25.7
$ws.Range("Q9").Value = 0.73
$ws.Range("R9").Value = 12.5
$ws.Range("S9").Value = 33.5
$ws.Range("T9").Value = 45.9
$ws.Range("U9").Value = 22.1
$ws.Range("V9").Value = 14.8
$ws.Range("W9").Value = 7.2
$ws.Range("Z9").Value = 22.5
$ws.Range("AB9").Value = 103.3
$ws.Range("AC9").Value = 0.6
$ws.Range("AD9").Value = 28
$ws.Range("AF9").Value = 12
$ws.Range("AG9").Value = 14
$ws.Range("AI9").Value = 13
$ws.Range("AP9").Value = 4
$ws.Range("AQ9").Value = 26
$ws.Range("AV9").Value = 12
$ws.Range("AW9").Value = 22
$ws.Range("AY9").Value = 20
$ws.Range("AZ9").Value = 27
$ws.Range("AD10").Value = 10
$ws.Range("AP10").Value = 6
$ws.Range("AH11").Value = 18
$ws.Range("AQ11").Value = 24
$ws.Range("BA11").Value = 18
$ws.Range("AZ12").Value = 12
$ws.Range("AI14").Value = 12
$ws.Range("AQ14").Value = 25
$ws.Range("AS14").Value = 12
$ws.Range("D15").Value = 42
$ws.Range("F15").Value = 26
$ws.Range("G15").Value = 0.381
$ws.Range("J15").Value = 83.7
$ws.Range("K15").Value = 0.44
$ws.Range("N15").Value = 0.368
$ws.Range("O15").Value = 17.7
$ws.Range("Q15").Value = 0.759
$ws.Range("R15").Value = 9.9
$ws.Range("S15").Value = 33
$ws.Range("T15").Value = 42.9
$ws.Range("V15").Value = 15.7
$ws.Range("W15").Value = 6.4
$ws.Range("X15").Value = 5.8
$ws.Range("AA15").Value = 19.3
$ws.Range("AD15").Value = 10
$ws.Range("AF15").Value = 22
$ws.Range("AH15").Value = 26
$ws.Range("AO15").Value = 14
$ws.Range("AQ15").Value = 15
$ws.Range("AS15").Value = 11
$ws.Range("AT15").Value = 17
$ws.Range("AU15").Value = 10
$ws.Range("AY15").Value = 15
$ws.Range("AZ15").Value = 13
$ws.Range("AD16").Value = 28
$ws.Range("D17").Value = 42
$ws.Range("E17").Value = 30
$ws.Range("G17").Value = 0.714
$ws.Range("I17").Value = 38.9
$ws.Range("K17").Value = 0.506
$ws.Range("O17").Value = 18.1
$ws.Range("P17").Value = 23.9
$ws.Range("Q17").Value = 0.759
$ws.Range("R17").Value = 7.1
$ws.Range("S17").Value = 29.5
$ws.Range("T17").Value = 36.5
$ws.Range("W17").Value = 9.199999999999999
$ws.Range("Y17").Value = 3.1
$ws.Range("Z17").Value = 20
$ws.Range("AA17").Value = 21.1
$ws.Range("AB17").Value = 104
$ws.Range("AC17").Value = 5.5
$ws.Range("AD17").Value = 10
$ws.Range("AI17").Value = 6
$ws.Range("AO17").Value = 11
$ws.Range("AQ17").Value = 14
$ws.Range("AS17").Value = 30
$ws.Range("AV17").Value = 16
$ws.Range("AZ17").Value = 8
$ws.Range("AQ18").Value = 19
$ws.Range("AF19").Value = 16
$ws.Range("AG19").Value = 16
$ws.Range("AO20").Value = 15
$ws.Range("AT20").Value = 19
$ws.Range("AV20").Value = 4
$ws.Range("BA20").Value = 17
$ws.Range("AD21").Value = 10
$ws.Range("AF21").Value = 24
$ws.Range("AQ21").Value = 16
$ws.Range("AS21").Value = 29
$ws.Range("AZ21").Value = 29
$ws.Range("AT22").Value = 2
$ws.Range("AQ23").Value = 18
$ws.Range("AD24").Value = 10
$ws.Range("AV25").Value = 14
$ws.Range("D26").Value = 42
$ws.Range("E26").Value = 31
$ws.Range("G26").Value = 0.738
$ws.Range("H26").Value = 48.4
$ws.Range("I26").Value = 40.5
$ws.Range("K26").Value = 0.458
$ws.Range("L26").Value = 10.1
$ws.Range("M26").Value = 25.6
$ws.Range("N26").Value = 0.395
$ws.Range("O26").Value = 18.2
$ws.Range("P26").Value = 22.2
$ws.Range("Q26").Value = 0.822
$ws.Range("R26").Value = 13
$ws.Range("S26").Value = 33.6
$ws.Range("T26").Value = 46.6
$ws.Range("U26").Value = 24.5
$ws.Range("V26").Value = 13.7
$ws.Range("AA26").Value = 20
$ws.Range("AD26").Value = 10
$ws.Range("AE26").Value = 4
$ws.Range("AH26").Value = 16
$ws.Range("AM26").Value = 3
$ws.Range("AO26").Value = 10
$ws.Range("AT26").Value = 1
$ws.Range("AV26").Value = 5
$ws.Range("BA26").Value = 19
$ws.Range("AP27").Value = 5
$ws.Range("AV27").Value = 15
$ws.Range("AD28").Value = 10
$ws.Range("AW29").Value = 23
$ws.Range("AZ29").Value = 28
$ws.Range("AQ30").Value = 17
$ws.Range("AV30").Value = 17
$ws.Range("AF31").Value = 16
$ws.Range("AG31").Value = 16
$ws.Range("AQ31").Value = 23
$ws.Range("AT31").Value = 21
